$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the 2021 column (R) one year further by adding a 2022 column (S),
# cloning R2:R5's formatting into S2:S5 so the new cells keep the same
# borders / number-formats / fonts as the existing year columns.
$ws.Range("R2:R5").Copy() | Out-Null
$ws.Range("S2:S5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data for 2022.
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 211650
$ws.Range("S5").Value = 2.9794303052841493

# Match the author's active selection on the new column.
$ws.Range("S2").Select() | Out-Null

Write-Host "done"
